$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the cell values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the target formatting (bold font, thin box border, center/top
# alignment) once on a scratch cell, then copy only the formatting onto
# each target cell. Doing it this way (vs. setting each property
# directly on every target cell) keeps the generated style table
# minimal/clean, matching a single shared cell style.
$scratch = $ws.Range("D10")
$scratch.Font.Bold = $true
$scratch.HorizontalAlignment = -4108  # xlCenter
$scratch.VerticalAlignment = -4160    # xlTop
$scratch.Borders.Weight = 2           # xlThin

$scratch.Copy()
foreach ($addr in @("B1", "A2", "B2")) {
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}
$scratch.Clear()
